$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text (e.g. "1.00", "69.234.55")
# Force the column to Text format first so Excel does not auto-convert
# numeric-looking strings into actual numbers, then restore the original
# (default/"Normal") style so no cell formatting is changed.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.234.55'
$ws.Range("D3").Value = '3.741.30'
$ws.Range("D5").Value = '602.52'
$ws.Range("D6").Value = '168.18'
$ws.Range("D7").Value = '3.739.54'
$ws.Range("D11").Value = '6.42'
$ws.Range("D13").Value = '38.04'
$ws.Range("D14").Value = '0.0000247'
$ws.Range("D15").Value = '4.365.02'
$ws.Range("D16").Value = '3.739.69'
$ws.Range("D17").Value = '69.199.58'
$ws.Range("D20").Value = '17.00'
$ws.Range("D22").Value = '494.92'
$ws.Range("D32").Value = '8.10'
$ws.Range("D33").Value = '31.62'
$ws.Range("D34").Value = '3.887.14'
$ws.Range("D36").Value = '3.673.90'
$ws.Range("D37").Value = '0.999'
$ws.Range("D38").Value = '1.02'
$ws.Range("D41").Value = '0.324'
$ws.Range("D42").Value = '3.00'
$ws.Range("D43").Value = '434.19'
$ws.Range("D44").Value = '48.67'
$ws.Range("D48").Value = '40.37'
$ws.Range("D49").Value = '140.85'
$ws.Range("D50").Value = '2.775.06'
$ws.Range("D51").Value = '0.0353'

$dRange.Style = "Normal"

# Volume(1h) (column E) percentage strings
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  +3.73%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  -3.56%  '
$ws.Range("E21").Value = '  +17.23%  '
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("E24").Value = '  +5.37%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +2.30%  '
$ws.Range("E31").Value = '  +6.37%  '
$ws.Range("E32").Value = '  +5.22%  '
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  +1.99%  '
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("E40").Value = '  +1.36%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  +5.98%  '
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("E51").Value = '  +0.75%  '

